$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.063.91"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "1.567.58"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  +0.73%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("E7").Value = "  +0.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("E12").Value = "  +1.09%  "

$ws.Range("D13").Value = "1.569.30"
$ws.Range("E13").Value = "  +1.28%  "

$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").Value = "27.056.77"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").Value = "0.0₃0704"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "

$ws.Range("E22").Value = "  +1.98%  "

$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("E29").Value = "  +0.69%  "

$ws.Range("E30").Value = "  +5.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.93%  "

$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("E33").Value = "  +2.95%  "

$ws.Range("D34").Value = "1.428.51"
$ws.Range("E34").Value = "  +1.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.47%  "

$ws.Range("E36").Value = "  +0.89%  "

$ws.Range("E37").Value = "  +3.42%  "

$ws.Range("E38").Value = "  +1.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.532"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("E40").Value = "  +2.79%  "

$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.06%  "

$ws.Range("E43").Value = "  +0.67%  "

$ws.Range("E44").Value = "  +0.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.28%  "

$ws.Range("E46").Value = "  +0.28%  "

$ws.Range("D47").Value = "1.706.02"
$ws.Range("E47").Value = "  +1.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("E49").Value = "  +1.66%  "

$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("E51").Value = "  +0.53%  "
